# New crime data collected — weekly CompStat refresh (Central Park Precinct)
# Moves the report forward one week: Volume 32 Number 6 -> Number 7,
# and the covered week from 2/3/2025-2/9/2025 to 2/10/2025-2/16/2025.
# Also refreshes the underlying crime-count table for rows 16-21 and 27-28.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text (rich-text shared strings collapse to plain strings here, but
# all runs already share one uniform font so nothing visually changes).
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 32   Number  7"
$ws.Range("C9").Value = "Report Covering the Week  2/10/2025  Through  2/16/2025"

# ---------------------------------------------------------------------------
# Helper donor cells that keep their original style/shared-string content
# throughout this script, so we can clone formatting (and, for the "0" /
# "***.*" placeholder text, the exact shared-string value) onto cells that
# need to flip between numeric and text representations.
#   C14 -> text "0"      (style 13)
#   E14 -> text "***.*"  (style 13)
#   F24 -> plain number  (style 14)
#   L24 -> percent number(style 15)
# ---------------------------------------------------------------------------
$xlPasteFormats = -4122

function Set-AsText($cellRef, $text) {
    $ws.Range("C14").Copy($ws.Range($cellRef))
    if ($text -eq "***.*") {
        $ws.Range("E14").Copy($ws.Range($cellRef))
    }
}

function Set-NumberStyle($cellRef, $value) {
    $ws.Range("F24").Copy()
    $ws.Range($cellRef).PasteSpecial($xlPasteFormats)
    $ws.Range($cellRef).Value = $value
}

function Set-PercentStyle($cellRef, $value) {
    $ws.Range("L24").Copy()
    $ws.Range($cellRef).PasteSpecial($xlPasteFormats)
    $ws.Range($cellRef).Value = $value
}

# ---------------------------------------------------------------------------
# Row 16 - Robbery: Week-to-date 2024 count/%chg go from real numbers to
# "no data" placeholders.
# ---------------------------------------------------------------------------
Set-AsText  "D16" "0"
Set-AsText  "E16" "***.*"

# ---------------------------------------------------------------------------
# Row 17 - Fel. Assault: placeholders become real counts for Week/28-Day/YTD.
# ---------------------------------------------------------------------------
Set-NumberStyle  "D17" 1
Set-PercentStyle "E17" -100
Set-NumberStyle  "G17" 1
Set-PercentStyle "H17" 100
Set-NumberStyle  "J17" 1
Set-PercentStyle "K17" 100
$ws.Range("N17").Value = -71.428571428571

# ---------------------------------------------------------------------------
# Row 19 - Gr. Larceny: new weekly counts; week-to-date 2024 flips to a
# placeholder.
# ---------------------------------------------------------------------------
$ws.Range("C19").Value = 1
Set-AsText "D19" "0"
Set-AsText "E19" "***.*"
$ws.Range("F19").Value = 4
$ws.Range("G19").Value = 1
$ws.Range("H19").Value = 300
$ws.Range("I19").Value = 4
$ws.Range("J19").Value = 2
$ws.Range("K19").Value = 100
$ws.Range("L19").Value = 33.333333333333
$ws.Range("M19").Value = 0
$ws.Range("N19").Value = -50

# ---------------------------------------------------------------------------
# Row 20 - G.L.A.: 32-year %chg placeholder becomes a real percentage.
# ---------------------------------------------------------------------------
Set-PercentStyle "N20" -100

# ---------------------------------------------------------------------------
# Row 21 - TOTAL: refreshed counts/percentages across the board.
# ---------------------------------------------------------------------------
$ws.Range("C21").Value = 1
$ws.Range("D21").Value = 1
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 6
$ws.Range("G21").Value = 5
$ws.Range("H21").Value = 20
$ws.Range("I21").Value = 6
$ws.Range("J21").Value = 7
$ws.Range("K21").Value = -14.285714285714
$ws.Range("L21").Value = 50
$ws.Range("M21").Value = 20
$ws.Range("N21").Value = -81.818181818181

# ---------------------------------------------------------------------------
# Row 27 - UCR Rape*: 28-day 2024 count/%chg revert to placeholders.
# ---------------------------------------------------------------------------
Set-AsText "G27" "0"
Set-AsText "H27" "***.*"

# ---------------------------------------------------------------------------
# Row 28 - Other Sex Crimes: week-to-date 2025 count reverts to a placeholder.
# ---------------------------------------------------------------------------
Set-AsText "C28" "0"
